# This script reshapes the "target" worksheet so that it stores a single,
# flat list of e-mail addresses (one header row + data) instead of the
# previous multi-column name/e-mail layout, and drops the now-unused
# wrap-text formatting that was only used by the removed rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old "TONAME" helper columns (C and G) and all of the old
#     e-mail values in column H except the header, together with any
#     formatting (e.g. the wrap-text style applied to the old H9 cell). ---
$ws.Range("C2:C8").Clear() | Out-Null
$ws.Range("G2:G8").Clear() | Out-Null
$ws.Range("H5:H28").Clear() | Out-Null

# --- Write the new, smaller set of e-mail addresses. ---
# First three addresses stay close to the header, in column H.
$ws.Range("H2").Value = "_anders@live.se"
$ws.Range("H3").Value = "_anton@live.se"
$ws.Range("H4").Value = "_armend_@live.se"

# Remaining addresses are written far down a new column (I), matching the
# large/sparse layout produced by the updated loading logic.
$ws.Range("I70259").Value = "kennard@gmail.com"
$ws.Range("I114093").Value = "watts@gmail.com"
$ws.Range("I145252").Value = "designs@gmail.com"
$ws.Range("I183666").Value = "michael.sun@gmail.com"
$ws.Range("I248024").Value = "samsalau@gmail.com"

# Restore the selection to F1 (previously F2), reflecting the header fix.
$ws.Range("F1").Select() | Out-Null
